$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5735
$ws.Range("F5").Value = 64
$ws.Range("F6").Value = 57
$ws.Range("F9").Value = 1551
$ws.Range("F10").Value = 12
$ws.Range("F12").Value = 654
$ws.Range("F13").Value = 1553
$ws.Range("F14").Value = 1553
$ws.Range("F15").Value = 1448
$ws.Range("F16").Value = 535
$ws.Range("F17").Value = 41
$ws.Range("F18").Value = 578
$ws.Range("F19").Value = 4137
$ws.Range("F20").Value = 4137
$ws.Range("F21").Value = 665
$ws.Range("F22").Value = 3314
$ws.Range("F23").Value = 786
$ws.Range("F25").Value = 2248
$ws.Range("F26").Value = 37
$ws.Range("F27").Value = 322
$ws.Range("F29").Value = 439
$ws.Range("F30").Value = 1199
$ws.Range("F31").Value = 775
$ws.Range("F33").Value = 1113
$ws.Range("F34").Value = 1130
$ws.Range("F35").Value = 79

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 101
$ws.Range("F17").Value = 117
$ws.Range("F19").Value = 213
$ws.Range("F20").Value = 482

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 573
$ws.Range("F5").Value = 112
$ws.Range("F6").Value = 204

$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 573
$ws.Range("F7").Value = 112
$ws.Range("F8").Value = 5735
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 57
$ws.Range("F14").Value = 101
$ws.Range("F19").Value = 1551
$ws.Range("F21").Value = 12
$ws.Range("F23").Value = 1553
$ws.Range("F24").Value = 1553
$ws.Range("F26").Value = 1448
$ws.Range("F27").Value = 535
$ws.Range("F28").Value = 41
$ws.Range("F29").Value = 578
$ws.Range("F31").Value = 4137
$ws.Range("F32").Value = 4137
$ws.Range("F33").Value = 665
$ws.Range("F34").Value = 3314
$ws.Range("F35").Value = 786
$ws.Range("F37").Value = 2248
$ws.Range("F38").Value = 37
$ws.Range("F40").Value = 439
$ws.Range("F41").Value = 1199
$ws.Range("F43").Value = 117
$ws.Range("F45").Value = 213
$ws.Range("F46").Value = 482
$ws.Range("F47").Value = 775
$ws.Range("F49").Value = 1113
$ws.Range("F50").Value = 1130
$ws.Range("F51").Value = 79
